# Applies the partner_466/txn_29/level_3 edit:
#  - Rename header C1 "Prin ID" -> "Partner ID"
#  - Rename header D1 "Prin Description" -> "Partner Description"
#  - Wrap several numeric-looking identifier columns (H,I,J,K,L,M,N,O) in
#    literal single quotes for rows 2-4 (e.g. 0010 -> '0010')
#  - Resize column D to fit its (unchanged) content
#  - Move the sheet selection from D2:D4 to M16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-write identifier columns as literal quoted text --------------
# A plain Value/Formula assignment of a string that *starts* with an
# apostrophe is interpreted by Excel as a quote-prefix (the leading quote
# is stripped and the cell is flagged as quote-prefixed text), which is not
# what the target file contains: it stores the literal characters
# "'0010'" (both quotes) as ordinary text with no quote-prefix flag.
# Building the string with CHAR(39) and then copying it as a *value* (not
# as typed user input) avoids that special-casing.
$staging = $ws.Range("Z1:Z4")
$staging.Cells.Item(1,1).Formula = "=CHAR(39)&""0010""&CHAR(39)"
$staging.Cells.Item(2,1).Formula = "=CHAR(39)&""0000""&CHAR(39)"
$staging.Cells.Item(3,1).Formula = "=CHAR(39)&""253""&CHAR(39)"
$staging.Cells.Item(4,1).Formula = "=CHAR(39)&""5511""&CHAR(39)"

$q0010 = $ws.Range("Z1")
$q0000 = $ws.Range("Z2")
$q253  = $ws.Range("Z3")
$q5511 = $ws.Range("Z4")

$targets0010 = @("H2","H3","H4","N2","N3","N4")
$targets0000 = @("I2","I3","I4","J2","J3","J4","M2","M3","M4","O2","O3","O4")
$targets253  = @("K2","K3","K4")
$targets5511 = @("L2","L3","L4")

$q0010.Copy()
foreach ($addr in $targets0010) {
    $ws.Range($addr).PasteSpecial(-4163, 0)
}

$q0000.Copy()
foreach ($addr in $targets0000) {
    $ws.Range($addr).PasteSpecial(-4163, 0)
}

$q253.Copy()
foreach ($addr in $targets253) {
    $ws.Range($addr).PasteSpecial(-4163, 0)
}

$q5511.Copy()
foreach ($addr in $targets5511) {
    $ws.Range($addr).PasteSpecial(-4163, 0)
}

$excel.CutCopyMode = 0
$staging.ClearContents()

# --- 2. Rename the two headers ------------------------------------------
$ws.Range("C1").Value = "Partner ID"
$ws.Range("D1").Value = "Partner Description"

# --- 3. Resize column D to fit the (unchanged) header/content -----------
$ws.Columns.Item(4).ColumnWidth = 16

# --- 4. Move the active selection ----------------------------------------
$null = $ws.Range("M16").Select()
